$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("A8").Value = 101
$ws.Range("B8").Value = "Australian Gouldian"
$ws.Range("C8").Value = "Coastal Cities"
$ws.Range("D8").Value = "12A"
$ws.Range("E8").Value = "Female"
$ws.Range("F8").Value = 223
$ws.Range("G8").Value = 111
$ws.Range("H8").Value = "15/05/2023"
$ws.Range("I8").Value = "Black"
$ws.Range("J8").Value = "Purple"
$ws.Range("K8").Value = "Green"
